# Add a new recipe row ("Base Yaourt Brassé Vache Sucré") to the
# ProductionRecipe sheet, mirroring the existing "Nature" recipe row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionRecipe")

$ws.Range("A3").Value = "REC_VACHE_BRASSE_SUCRE"
$ws.Range("B3").Value = "Base Yaourt Brassé Vache Sucré"
$ws.Range("C3").Value = "BASE_VACHE_BRASSE_SUCRE"
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Recette de la Base Blanche de Yaourt Brassé Vache Sucré"

# Match the author's final selection state (A1:F3 highlighted).
$ws.Range("A1:F3").Select()
